$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: new "git add" note (set first so it reuses the freed shared-string slot) ---
$ws2.Range("T22").Value = "git add"

# --- Sheet1: "skeletton environment" split into 3 parts + new setRoll task ---

# Row 11
$ws1.Range("C11").Value = 0.28472222222222221
$ws1.Range("C11").NumberFormat = $ws1.Range("D6").NumberFormat
$ws1.Range("F11").Value = 23
$ws1.Range("I11").Value = "create skeletton environment part 1"

# Row 12
$ws1.Range("B12").Formula = "=SUM(F6:F25)"
$ws1.Range("C12").Value = 0.33888888888888885
$ws1.Range("C12").NumberFormat = $ws1.Range("D6").NumberFormat
$ws1.Range("D12").Value = 0.29166666666666669
$ws1.Range("D12").NumberFormat = $ws1.Range("D6").NumberFormat
$ws1.Range("F12").Value = 68
$ws1.Range("I12").Value = "create skeletton environment part 2"

# Row 13
$ws1.Range("B13").Formula = "=B12/60"
$ws1.Range("C13").NumberFormat = $ws1.Range("D6").NumberFormat
$ws1.Range("D13").NumberFormat = $ws1.Range("D6").NumberFormat
$ws1.Range("F13").Formula = "=15+18+41+15+60"
$ws1.Range("I13").Value = "create skeletton environment part 3"

# Row 14
$ws1.Range("F14").Value = 30
$ws1.Range("I14").Value = "Implement Function setRoll <Frame>"

# Row 15: task label cleared (text moved further down the list)
$ws1.Range("I15").Value = ""

# Rows 16-19: remaining task labels shift down by four rows
$ws1.Range("I16").Value = "Create function to calcule score of a given game "
$ws1.Range("I17").Value = "create setRoll function "
$ws1.Range("I18").Value = "create function that check score "
$ws1.Range("I19").Value = "create function that print frames and final score and score after every frame"

# Rows 20-23: newly appended task labels
$ws1.Range("I20").Value = "add multiple player feature "
$ws1.Range("I21").Value = "add rank after gameover"
$ws1.Range("I22").Value = "calcule score after every roll setted in the game "
$ws1.Range("I23").Value = "next feature pins presentation feature "

# --- View state: Sheet1 becomes the active/visible tab ---
$ws1.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws1.Range("F14").Select() | Out-Null
